$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsCurrency = $wb.Worksheets.Item("Currency Conversion")

# --- Currency Conversion sheet: add 2018 USD/INR conversion rate row ---
$wsCurrency.Range("A12").Value = 2018
$wsCurrency.Range("B12").NumberFormat = "0.00"
$wsCurrency.Range("B12").Value = 68.657300000000006

# --- Currency Conversion sheet: add 2018 inflation rate row ---
$wsCurrency.Range("A36").Value = 2018
$wsCurrency.Range("B36").NumberFormat = "0.00%"
$wsCurrency.Range("B36").Value = 0.0524
$wsCurrency.Range("C36").Formula = "=C35*(1+B36)"

# --- About sheet: turn the RBI reference-rate text (B7) into a hyperlink ---
[void]$wsAbout.Hyperlinks.Add($wsAbout.Range("B7"), "https://rbi.org.in/scripts/BS_DisplayReferenceRate.aspx")
$wsAbout.Range("B7").Style = "Hyperlink"

# --- Restore sheet selections / active sheet to match the author's last save ---
[void]$wsAbout.Range("B8").Select()
[void]$wsCurrency.Activate()
[void]$wsCurrency.Range("E11").Select()
